# adicionando uma interface gráfica simples
#
# Reverts the worksheet's header/data styling back to the plain
# "bold header + thin border" look (dropping the fancy dark-blue /
# light-blue theme that had been applied), removes the second sample
# row, renames the first sample file, and stores the OAB number as
# plain text instead of a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Drop the "processo_4" sample row (row 3) entirely.
# ------------------------------------------------------------------
$ws.Rows.Item(3).Delete()

# ------------------------------------------------------------------
# 2) Strip the fancy dark-blue header fill/font back down to a plain
#    bold font with a thin border (center/top aligned).
# ------------------------------------------------------------------
$header = $ws.Range("A1:F1")
$header.ClearFormats()
$header.Font.Bold = $true
$header.Borders.LineStyle = 1
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160

# ------------------------------------------------------------------
# 3) Strip the light-blue fill from the remaining data row so it goes
#    back to the workbook's default (unstyled) look.
# ------------------------------------------------------------------
$data = $ws.Range("A2:F2")
$data.ClearFormats()

# ------------------------------------------------------------------
# 4) Rename the sample file and store OAB as text ("44432") rather
#    than a number.
# ------------------------------------------------------------------
$ws.Range("A2").Value = "copy_processo_1"

$oab = $ws.Range("E2")
$oab.NumberFormat = "@"
$oab.Value = "44432"
$oab.ClearFormats()
